$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
